$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new entry "김일현" to the next empty row (A3), matching the reverted
# commit's original data: A1="김주현", A2="김순주", A3="김일현"
$ws.Range("A3").Value = "김일현"

# Move the active cell selection down to A4, as in the original file
$ws.Range("A4").Select()
